{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Load the cells we need to inspect/modify.\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst headerCols = 5;\nconst dataRows = table.rowCount; // includes header row at index 0\n\n// 1) Header row: italic column-title cells shrink from 11pt (sz 22) to 10pt (sz 20).\nfor (let c = 0; c < headerCols; c++) {\n  const cell = table.getCell(0, c);\n  const paras = cell.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n  para.getRange(\"Content\").font.size = 10;\n}\nawait context.sync();\n\n// 2) Data rows: first column (\"STOCK\" values) right-align instead of left-align.\nfor (let r = 1; r < dataRows; r++) {\n  const cell = table.getCell(r, 0);\n  const paras = cell.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const para = paras.items[0];\n  para.alignment = Word.Alignment.right;\n}\nawait context.sync();\n\n// 3) Pad specific numeric cell values to 4 decimal places.\nconst fixups = [\n  { row: 2, col: 1, text: \"0.0000\" },\n  { row: 2, col: 3, text: \"0.5000\" },\n  { row: 4, col: 3, text: \"1.0000\" },\n  { row: 7, col: 1, text: \"1.0000\" },\n  { row: 8, col: 1, text: \"0.0440\" },\n  { row: 9, col: 4, text: \"0.5040\" },\n  { row: 10, col: 3, text: \"0.0000\" },\n  { row: 11, col: 4, text: \"0.5040\" }\n];\n\nfor (const fx of fixups) {\n  const cell = table.getCell(fx.row, fx.col);\n  cell.value = fx.text;\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1) Header row: italic column-title cells shrink from 11pt (sz 22) to 10pt (sz 20)\nfor ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Cell(1, $c).Range.Font.Size = 10\n}\n\n# 2) Data rows (2..12): first column (\"STOCK\" values) right-align instead of left-align\nfor ($r = 2; $r -le $t.Rows.Count; $r++) {\n    $t.Cell($r, 1).Range.ParagraphFormat.Alignment = 2\n}\n\n# 3) Pad specific numeric cell values to 4 decimal places\n$fixups = @(\n    @{ Row = 3;  Col = 2; Text = \"0.0000\" },\n    @{ Row = 3;  Col = 4; Text = \"0.5000\" },\n    @{ Row = 5;  Col = 4; Text = \"1.0000\" },\n    @{ Row = 8;  Col = 2; Text = \"1.0000\" },\n    @{ Row = 9;  Col = 2; Text = \"0.0440\" },\n    @{ Row = 10; Col = 5; Text = \"0.5040\" },\n    @{ Row = 11; Col = 4; Text = \"0.0000\" },\n    @{ Row = 12; Col = 5; Text = \"0.5040\" }\n)\n\nforeach ($fx in $fixups) {\n    $t.Cell($fx.Row, $fx.Col).Range.Text = $fx.Text\n}\n"}
